# Insert two new rows at row 624 (pushing the existing rows 624-717 down to
# become rows 626-719), then populate the two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 624.
$ws.Range("624:625").Insert()

# --- New row 624 ---
$ws.Cells.Item(624, 1).Value  = 10
$ws.Cells.Item(624, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(624, 3).Value  = "La Araucanía"
$ws.Cells.Item(624, 4).Value  = 44522
$ws.Cells.Item(624, 5).Value  = 9
$ws.Cells.Item(624, 6).Value  = 100112004
$ws.Cells.Item(624, 7).Value  = "Cebolla"
$ws.Cells.Item(624, 8).Value  = "Sin especificar"
$ws.Cells.Item(624, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(624, 10).Value = 600
$ws.Cells.Item(624, 11).Value = 5000
$ws.Cells.Item(624, 12).Value = 5000
$ws.Cells.Item(624, 13).Value = 5000
$ws.Cells.Item(624, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(624, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(624, 16).Value = 278
$ws.Cells.Item(624, 17).Value = 18
$ws.Cells.Item(624, 18).Value = "Hortaliza"

# --- New row 625 ---
$ws.Cells.Item(625, 1).Value  = 10
$ws.Cells.Item(625, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(625, 3).Value  = "La Araucanía"
$ws.Cells.Item(625, 4).Value  = 44522
$ws.Cells.Item(625, 5).Value  = 9
$ws.Cells.Item(625, 6).Value  = 100112004
$ws.Cells.Item(625, 7).Value  = "Cebolla"
$ws.Cells.Item(625, 8).Value  = "Sin especificar"
$ws.Cells.Item(625, 9).Value  = "Primera"
$ws.Cells.Item(625, 10).Value = 1000
$ws.Cells.Item(625, 11).Value = 4500
$ws.Cells.Item(625, 12).Value = 4500
$ws.Cells.Item(625, 13).Value = 4500
$ws.Cells.Item(625, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(625, 15).Value = "Perú"
$ws.Cells.Item(625, 16).Value = 250
$ws.Cells.Item(625, 17).Value = 18
$ws.Cells.Item(625, 18).Value = "Hortaliza"

Write-Output "Inserted rows 624-625; new dimension should be A1:R719"
